$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data to the latest snapshot.
# Column D (Price) values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (preserving formats like "0.130" / "43.383.67")
# instead of auto-converting them into floating point numbers.

$ws.Range("D2").Value = "'43.383.67"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "'2.243.14"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'230.52"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'0.643"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("D7").Value = "'63.81"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "'0.0955"
$ws.Range("E10").Value = "  -5.46%  "
$ws.Range("D11").Value = "'56.42"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "'27.32"
$ws.Range("E12").Value = "  +5.32%  "
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Value = "'2.575.25"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "'15.24"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "'0.825"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'2.241.98"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "'43.296.82"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").Value = "'0.0₃0964"
$ws.Range("E20").Value = "  -5.56%  "
$ws.Range("D21").Value = "'72.98"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'246.27"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'3.78"
$ws.Range("E25").Value = "  +33.51%  "
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("E27").Value = "  -4.49%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.71"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'173.88"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "'21.66"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("D31").Value = "'0.130"
$ws.Range("E31").Value = "  -5.29%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'4.91"
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").Value = "'4.91"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -7.46%  "
$ws.Range("D38").Value = "'6.31"
$ws.Range("E38").Value = "  -5.87%  "
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("E43").Value = "  +3.66%  "
$ws.Range("D44").Value = "'17.02"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").Value = "'96.53"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "'1.18"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("B48").Value = "TerraClassic"
$ws.Range("C48").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D48").Value = "'0.000208"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").Value = "'9.93"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'1.430.95"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "'2.25"
$ws.Range("E51").Value = "  -2.53%  "
